# ---- Build "2022-Q1" sheet with fund holdings data ----
$wb = $excel.ActiveWorkbook
$sumSheet = $wb.Worksheets.Item(6)
$ws = $wb.Worksheets.Add($sumSheet)
$ws.Name = "2022-Q1"

$src = $wb.Worksheets.Item(5)
$src.Range("A1:H26").Copy()
$ws.Range("A1:H26").PasteSpecial(-4122)
$ws.Cells.Item(1,1).Clear()

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Header row
$ws.Cells.Item(1,2).Value = '基金代码'
$ws.Cells.Item(1,3).Value = '基金名称'
$ws.Cells.Item(1,4).Value = '基金规模'
$ws.Cells.Item(1,5).Value = '股票总仓位'
$ws.Cells.Item(1,6).Value = '仓位占比'
$ws.Cells.Item(1,7).Value = '持有市值(亿元)'
$ws.Cells.Item(1,8).Value = '仓位排名'

# Row 2: 006408 汇添富消费升级混合
$ws.Cells.Item(2,1).Value = 0
Set-TextValue $ws.Cells.Item(2,2) '006408'
$ws.Cells.Item(2,3).Value = '汇添富消费升级混合'
Set-TextValue $ws.Cells.Item(2,4) '42.63'
Set-TextValue $ws.Cells.Item(2,5) '87.97'
Set-TextValue $ws.Cells.Item(2,6) '4.64'
Set-TextValue $ws.Cells.Item(2,7) '1.9780'
$ws.Cells.Item(2,8).Value = 6

# Row 3: 009086 鹏华价值共赢两年持有期混合
$ws.Cells.Item(3,1).Value = 1
Set-TextValue $ws.Cells.Item(3,2) '009086'
$ws.Cells.Item(3,3).Value = '鹏华价值共赢两年持有期混合'
Set-TextValue $ws.Cells.Item(3,4) '19.26'
Set-TextValue $ws.Cells.Item(3,5) '64.26'
Set-TextValue $ws.Cells.Item(3,6) '5.98'
Set-TextValue $ws.Cells.Item(3,7) '1.1517'
$ws.Cells.Item(3,8).Value = 2

# Row 4: 660010 农银策略精选混合
$ws.Cells.Item(4,1).Value = 2
Set-TextValue $ws.Cells.Item(4,2) '660010'
$ws.Cells.Item(4,3).Value = '农银策略精选混合'
Set-TextValue $ws.Cells.Item(4,4) '31.49'
Set-TextValue $ws.Cells.Item(4,5) '75.60'
Set-TextValue $ws.Cells.Item(4,6) '3.37'
Set-TextValue $ws.Cells.Item(4,7) '1.0612'
$ws.Cells.Item(4,8).Value = 10

# Row 5: 010815 农银汇理新兴消费股票
$ws.Cells.Item(5,1).Value = 3
Set-TextValue $ws.Cells.Item(5,2) '010815'
$ws.Cells.Item(5,3).Value = '农银汇理新兴消费股票'
Set-TextValue $ws.Cells.Item(5,4) '29.81'
Set-TextValue $ws.Cells.Item(5,5) '83.28'
Set-TextValue $ws.Cells.Item(5,6) '3.00'
Set-TextValue $ws.Cells.Item(5,7) '0.8943'
$ws.Cells.Item(5,8).Value = 10

# Row 6: 000127 农银行业领先混合
$ws.Cells.Item(6,1).Value = 4
Set-TextValue $ws.Cells.Item(6,2) '000127'
$ws.Cells.Item(6,3).Value = '农银行业领先混合'
Set-TextValue $ws.Cells.Item(6,4) '12.28'
Set-TextValue $ws.Cells.Item(6,5) '75.75'
Set-TextValue $ws.Cells.Item(6,6) '3.38'
Set-TextValue $ws.Cells.Item(6,7) '0.4151'
$ws.Cells.Item(6,8).Value = 10

# Row 7: 240001 华宝宝康消费品混合
$ws.Cells.Item(7,1).Value = 5
Set-TextValue $ws.Cells.Item(7,2) '240001'
$ws.Cells.Item(7,3).Value = '华宝宝康消费品混合'
Set-TextValue $ws.Cells.Item(7,4) '11.24'
Set-TextValue $ws.Cells.Item(7,5) '62.57'
Set-TextValue $ws.Cells.Item(7,6) '3.69'
Set-TextValue $ws.Cells.Item(7,7) '0.4148'
$ws.Cells.Item(7,8).Value = 4

# Row 8: 011153 华宝新兴消费混合A
$ws.Cells.Item(8,1).Value = 6
Set-TextValue $ws.Cells.Item(8,2) '011153'
$ws.Cells.Item(8,3).Value = '华宝新兴消费混合A'
Set-TextValue $ws.Cells.Item(8,4) '6.63'
Set-TextValue $ws.Cells.Item(8,5) '87.87'
Set-TextValue $ws.Cells.Item(8,6) '4.68'
Set-TextValue $ws.Cells.Item(8,7) '0.3103'
$ws.Cells.Item(8,8).Value = 3

# Row 9: 008819 农银汇理策略趋势混合
$ws.Cells.Item(9,1).Value = 7
Set-TextValue $ws.Cells.Item(9,2) '008819'
$ws.Cells.Item(9,3).Value = '农银汇理策略趋势混合'
Set-TextValue $ws.Cells.Item(9,4) '6.17'
Set-TextValue $ws.Cells.Item(9,5) '76.46'
Set-TextValue $ws.Cells.Item(9,6) '3.48'
Set-TextValue $ws.Cells.Item(9,7) '0.2147'
$ws.Cells.Item(9,8).Value = 9

# Row 10: 004634 新疆前海联合泳涛灵活配置混合A
$ws.Cells.Item(10,1).Value = 8
Set-TextValue $ws.Cells.Item(10,2) '004634'
$ws.Cells.Item(10,3).Value = '新疆前海联合泳涛灵活配置混合A'
Set-TextValue $ws.Cells.Item(10,4) '1.33'
Set-TextValue $ws.Cells.Item(10,5) '89.65'
Set-TextValue $ws.Cells.Item(10,6) '7.58'
Set-TextValue $ws.Cells.Item(10,7) '0.1008'
$ws.Cells.Item(10,8).Value = 3

# Row 11: 510630 华夏上证主要消费ETF
$ws.Cells.Item(11,1).Value = 9
Set-TextValue $ws.Cells.Item(11,2) '510630'
$ws.Cells.Item(11,3).Value = '华夏上证主要消费ETF'
Set-TextValue $ws.Cells.Item(11,4) '3.36'
Set-TextValue $ws.Cells.Item(11,5) '99.52'
Set-TextValue $ws.Cells.Item(11,6) '2.27'
Set-TextValue $ws.Cells.Item(11,7) '0.0763'
$ws.Cells.Item(11,8).Value = 10

# Row 12: 012080 易方达中证500指数量化增强型证券投资基金A
$ws.Cells.Item(12,1).Value = 10
Set-TextValue $ws.Cells.Item(12,2) '012080'
$ws.Cells.Item(12,3).Value = '易方达中证500指数量化增强型证券投资基金A'
Set-TextValue $ws.Cells.Item(12,4) '6.82'
Set-TextValue $ws.Cells.Item(12,5) '84.83'
Set-TextValue $ws.Cells.Item(12,6) '0.99'
Set-TextValue $ws.Cells.Item(12,7) '0.0675'
$ws.Cells.Item(12,8).Value = 4

# Row 13: 159758 华夏中证红利质量ETF
$ws.Cells.Item(13,1).Value = 11
Set-TextValue $ws.Cells.Item(13,2) '159758'
$ws.Cells.Item(13,3).Value = '华夏中证红利质量ETF'
Set-TextValue $ws.Cells.Item(13,4) '1.81'
Set-TextValue $ws.Cells.Item(13,5) '99.16'
Set-TextValue $ws.Cells.Item(13,6) '3.56'
Set-TextValue $ws.Cells.Item(13,7) '0.0644'
$ws.Cells.Item(13,8).Value = 7

# Row 14: 010636 财通安盈混合A
$ws.Cells.Item(14,1).Value = 12
Set-TextValue $ws.Cells.Item(14,2) '010636'
$ws.Cells.Item(14,3).Value = '财通安盈混合A'
Set-TextValue $ws.Cells.Item(14,4) '3.81'
Set-TextValue $ws.Cells.Item(14,5) '34.52'
Set-TextValue $ws.Cells.Item(14,6) '1.16'
Set-TextValue $ws.Cells.Item(14,7) '0.0442'
$ws.Cells.Item(14,8).Value = 9

# Row 15: 006235 东方城镇消费主题混合
$ws.Cells.Item(15,1).Value = 13
Set-TextValue $ws.Cells.Item(15,2) '006235'
$ws.Cells.Item(15,3).Value = '东方城镇消费主题混合'
Set-TextValue $ws.Cells.Item(15,4) '0.50'
Set-TextValue $ws.Cells.Item(15,5) '90.32'
Set-TextValue $ws.Cells.Item(15,6) '7.51'
Set-TextValue $ws.Cells.Item(15,7) '0.0376'
$ws.Cells.Item(15,8).Value = 2

# Row 16: 001421 南方量化成长股票
$ws.Cells.Item(16,1).Value = 14
Set-TextValue $ws.Cells.Item(16,2) '001421'
$ws.Cells.Item(16,3).Value = '南方量化成长股票'
Set-TextValue $ws.Cells.Item(16,4) '1.70'
Set-TextValue $ws.Cells.Item(16,5) '92.11'
Set-TextValue $ws.Cells.Item(16,6) '1.88'
Set-TextValue $ws.Cells.Item(16,7) '0.0320'
$ws.Cells.Item(16,8).Value = 5

# Row 17: 008353 泰达宏利消费行业量化精选混合A
$ws.Cells.Item(17,1).Value = 15
Set-TextValue $ws.Cells.Item(17,2) '008353'
$ws.Cells.Item(17,3).Value = '泰达宏利消费行业量化精选混合A'
Set-TextValue $ws.Cells.Item(17,4) '0.52'
Set-TextValue $ws.Cells.Item(17,5) '92.25'
Set-TextValue $ws.Cells.Item(17,6) '5.92'
Set-TextValue $ws.Cells.Item(17,7) '0.0308'
$ws.Cells.Item(17,8).Value = 2

# Row 18: 010637 财通安盈混合C
$ws.Cells.Item(18,1).Value = 16
Set-TextValue $ws.Cells.Item(18,2) '010637'
$ws.Cells.Item(18,3).Value = '财通安盈混合C'
Set-TextValue $ws.Cells.Item(18,4) '1.83'
Set-TextValue $ws.Cells.Item(18,5) '34.52'
Set-TextValue $ws.Cells.Item(18,6) '1.16'
Set-TextValue $ws.Cells.Item(18,7) '0.0212'
$ws.Cells.Item(18,8).Value = 9

# Row 19: 012081 易方达中证500指数量化增强型证券投资基金C
$ws.Cells.Item(19,1).Value = 17
Set-TextValue $ws.Cells.Item(19,2) '012081'
$ws.Cells.Item(19,3).Value = '易方达中证500指数量化增强型证券投资基金C'
Set-TextValue $ws.Cells.Item(19,4) '1.57'
Set-TextValue $ws.Cells.Item(19,5) '84.83'
Set-TextValue $ws.Cells.Item(19,6) '0.99'
Set-TextValue $ws.Cells.Item(19,7) '0.0155'
$ws.Cells.Item(19,8).Value = 4

# Row 20: 561350 国泰中证500ETF
$ws.Cells.Item(20,1).Value = 18
Set-TextValue $ws.Cells.Item(20,2) '561350'
$ws.Cells.Item(20,3).Value = '国泰中证500ETF'
Set-TextValue $ws.Cells.Item(20,4) '1.78'
Set-TextValue $ws.Cells.Item(20,5) '96.89'
Set-TextValue $ws.Cells.Item(20,6) '0.65'
Set-TextValue $ws.Cells.Item(20,7) '0.0116'
$ws.Cells.Item(20,8).Value = 6

# Row 21: 008354 泰达宏利消费行业量化精选混合C
$ws.Cells.Item(21,1).Value = 19
Set-TextValue $ws.Cells.Item(21,2) '008354'
$ws.Cells.Item(21,3).Value = '泰达宏利消费行业量化精选混合C'
Set-TextValue $ws.Cells.Item(21,4) '0.12'
Set-TextValue $ws.Cells.Item(21,5) '92.25'
Set-TextValue $ws.Cells.Item(21,6) '5.92'
Set-TextValue $ws.Cells.Item(21,7) '0.0071'
$ws.Cells.Item(21,8).Value = 2

# Row 22: 011154 华宝新兴消费混合C
$ws.Cells.Item(22,1).Value = 20
Set-TextValue $ws.Cells.Item(22,2) '011154'
$ws.Cells.Item(22,3).Value = '华宝新兴消费混合C'
Set-TextValue $ws.Cells.Item(22,4) '0.12'
Set-TextValue $ws.Cells.Item(22,5) '87.87'
Set-TextValue $ws.Cells.Item(22,6) '4.68'
Set-TextValue $ws.Cells.Item(22,7) '0.0056'
$ws.Cells.Item(22,8).Value = 3

# Row 23: 007943 富安达中证 500 指数增强
$ws.Cells.Item(23,1).Value = 21
Set-TextValue $ws.Cells.Item(23,2) '007943'
$ws.Cells.Item(23,3).Value = '富安达中证 500 指数增强'
Set-TextValue $ws.Cells.Item(23,4) '0.21'
Set-TextValue $ws.Cells.Item(23,5) '93.50'
Set-TextValue $ws.Cells.Item(23,6) '1.33'
Set-TextValue $ws.Cells.Item(23,7) '0.0028'
$ws.Cells.Item(23,8).Value = 7

# Row 24: 010957 九泰久安量化股票型证券投资基金A
$ws.Cells.Item(24,1).Value = 22
Set-TextValue $ws.Cells.Item(24,2) '010957'
$ws.Cells.Item(24,3).Value = '九泰久安量化股票型证券投资基金A'
Set-TextValue $ws.Cells.Item(24,4) '0.01'
Set-TextValue $ws.Cells.Item(24,5) '81.02'
Set-TextValue $ws.Cells.Item(24,6) '3.36'
Set-TextValue $ws.Cells.Item(24,7) '0.0003'
$ws.Cells.Item(24,8).Value = 10

# Row 25: 010961 九泰久安量化股票型证券投资基金C
$ws.Cells.Item(25,1).Value = 23
Set-TextValue $ws.Cells.Item(25,2) '010961'
$ws.Cells.Item(25,3).Value = '九泰久安量化股票型证券投资基金C'
Set-TextValue $ws.Cells.Item(25,4) '0.00'
Set-TextValue $ws.Cells.Item(25,5) '81.02'
Set-TextValue $ws.Cells.Item(25,6) '3.36'
$ws.Cells.Item(25,7).Value = 0
$ws.Cells.Item(25,8).Value = 10

# Row 26: 007041 新疆前海联合泳涛灵活配置混合C
$ws.Cells.Item(26,1).Value = 24
Set-TextValue $ws.Cells.Item(26,2) '007041'
$ws.Cells.Item(26,3).Value = '新疆前海联合泳涛灵活配置混合C'
Set-TextValue $ws.Cells.Item(26,4) '0.00'
Set-TextValue $ws.Cells.Item(26,5) '89.65'
Set-TextValue $ws.Cells.Item(26,6) '7.58'
$ws.Cells.Item(26,7).Value = 0
$ws.Cells.Item(26,8).Value = 3


# ---- Insert new row in "总计" sheet for 2022-Q1 summary ----
$totalSheet = $wb.Worksheets.Item(7)
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Cells.Item(3,1).Copy()
$totalSheet.Cells.Item(2,1).PasteSpecial(-4122)
$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 25
$totalSheet.Cells.Item(2,4).Value = 6.96

# Re-number the index column (A) for the rows that shifted down
$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(4,1).Value = 2
$totalSheet.Cells.Item(5,1).Value = 3
$totalSheet.Cells.Item(6,1).Value = 4
$totalSheet.Cells.Item(7,1).Value = 5

